$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1581976666666667
$ws.Range("H2").Value = 0.474593
$ws.Range("I2").Value = 0.1400666049254827
$ws.Range("J2").Value = 0.1400666049254826
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.949891
$ws.Range("N2").Value = 2.849673
$ws.Range("O2").Value = 0.1664393778377885
$ws.Range("P2").Value = 0.1664393778377885
$ws.Range("Q2").Value = 0.1502705397876667
$ws.Range("R2").Value = 1.352434858089
$ws.Range("S2").Value = 0.02331259857964866
$ws.Range("T2").Value = 0.02331259857964866

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1581976666666667
$ws.Range("H3").Value = 0.474593
$ws.Range("I3").Value = 0.1400666049254827
$ws.Range("J3").Value = 0.1400666049254826
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.834886333333333
$ws.Range("N3").Value = 8.504659
$ws.Range("O3").Value = 0.4967272219242518
$ws.Range("P3").Value = 0.4967272219242519
$ws.Range("Q3").Value = 0.4484724031985556
$ws.Range("R3").Value = 4.036251628787
$ws.Range("S3").Value = 0.06957489554899673
$ws.Range("T3").Value = 0.06957489554899673

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1581976666666667
$ws.Range("H4").Value = 0.474593
$ws.Range("I4").Value = 0.1400666049254827
$ws.Range("J4").Value = 0.1400666049254826
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.722217666666667
$ws.Range("N4").Value = 5.166653
$ws.Range("O4").Value = 0.3017660309880268
$ws.Range("P4").Value = 0.3017660309880269
$ws.Range("Q4").Value = 0.2724508163587778
$ws.Range("R4").Value = 2.452057347229
$ws.Range("S4").Value = 0.04226734344233092
$ws.Range("T4").Value = 0.04226734344233092

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1581976666666667
$ws.Range("H5").Value = 0.474593
$ws.Range("I5").Value = 0.1400666049254827
$ws.Range("J5").Value = 0.1400666049254826
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.200134
$ws.Range("N5").Value = 0.600402
$ws.Range("O5").Value = 0.03506736924993285
$ws.Range("P5").Value = 0.03506736924993285
$ws.Range("Q5").Value = 0.03166073182066667
$ws.Range("R5").Value = 0.284946586386
$ws.Range("S5").Value = 0.004911767354506364
$ws.Range("T5").Value = 0.004911767354506363

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7939349999999999
$ws.Range("H6").Value = 2.381805
$ws.Range("I6").Value = 0.7029419733214338
$ws.Range("J6").Value = 0.7029419733214337
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.949891
$ws.Range("N6").Value = 2.849673
$ws.Range("O6").Value = 0.1664393778377885
$ws.Range("P6").Value = 0.1664393778377885
$ws.Range("Q6").Value = 0.754151711085
$ws.Range("R6").Value = 6.787365399765
$ws.Range("S6").Value = 0.1169972246956868
$ws.Range("T6").Value = 0.1169972246956868

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7939349999999999
$ws.Range("H7").Value = 2.381805
$ws.Range("I7").Value = 0.7029419733214338
$ws.Range("J7").Value = 0.7029419733214337
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.834886333333333
$ws.Range("N7").Value = 8.504659
$ws.Range("O7").Value = 0.4967272219242518
$ws.Range("P7").Value = 0.4967272219242519
$ws.Range("Q7").Value = 2.250715481055
$ws.Range("R7").Value = 20.256439329495
$ws.Range("S7").Value = 0.3491704135819074
$ws.Range("T7").Value = 0.3491704135819074

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.7939349999999999
$ws.Range("H8").Value = 2.381805
$ws.Range("I8").Value = 0.7029419733214338
$ws.Range("J8").Value = 0.7029419733214337
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.722217666666667
$ws.Range("N8").Value = 5.166653
$ws.Range("O8").Value = 0.3017660309880268
$ws.Range("P8").Value = 0.3017660309880269
$ws.Range("Q8").Value = 1.367328883185
$ws.Range("R8").Value = 12.305959948665
$ws.Range("S8").Value = 0.2121240093041006
$ws.Range("T8").Value = 0.2121240093041006

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.7939349999999999
$ws.Range("H9").Value = 2.381805
$ws.Range("I9").Value = 0.7029419733214338
$ws.Range("J9").Value = 0.7029419733214337
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.200134
$ws.Range("N9").Value = 0.600402
$ws.Range("O9").Value = 0.03506736924993285
$ws.Range("P9").Value = 0.03506736924993285
$ws.Range("Q9").Value = 0.15889338729
$ws.Range("R9").Value = 1.43004048561
$ws.Range("S9").Value = 0.02465032573973916
$ws.Range("T9").Value = 0.02465032573973916

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1199896666666667
$ws.Range("H10").Value = 0.359969
$ws.Range("I10").Value = 0.106237630366274
$ws.Range("J10").Value = 0.106237630366274
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.949891
$ws.Range("N10").Value = 2.849673
$ws.Range("O10").Value = 0.1664393778377885
$ws.Range("P10").Value = 0.1664393778377885
$ws.Range("Q10").Value = 0.1139771044596667
$ws.Range("R10").Value = 1.025793940137
$ws.Range("S10").Value = 0.01768212510112359
$ws.Range("T10").Value = 0.01768212510112359

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1199896666666667
$ws.Range("H11").Value = 0.359969
$ws.Range("I11").Value = 0.106237630366274
$ws.Range("J11").Value = 0.106237630366274
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.834886333333333
$ws.Range("N11").Value = 8.504659
$ws.Range("O11").Value = 0.4967272219242518
$ws.Range("P11").Value = 0.4967272219242519
$ws.Range("Q11").Value = 0.3401570661745555
$ws.Range("R11").Value = 3.061413595571
$ws.Range("S11").Value = 0.05277112299565481
$ws.Range("T11").Value = 0.05277112299565481

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.1199896666666667
$ws.Range("H12").Value = 0.359969
$ws.Range("I12").Value = 0.106237630366274
$ws.Range("J12").Value = 0.106237630366274
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.722217666666667
$ws.Range("N12").Value = 5.166653
$ws.Range("O12").Value = 0.3017660309880268
$ws.Range("P12").Value = 0.3017660309880269
$ws.Range("Q12").Value = 0.2066483237507778
$ws.Range("R12").Value = 1.859834913757
$ws.Range("S12").Value = 0.03205890805720358
$ws.Range("T12").Value = 0.03205890805720357

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.1199896666666667
$ws.Range("H13").Value = 0.359969
$ws.Range("I13").Value = 0.106237630366274
$ws.Range("J13").Value = 0.106237630366274
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.200134
$ws.Range("N13").Value = 0.600402
$ws.Range("O13").Value = 0.03506736924993285
$ws.Range("P13").Value = 0.03506736924993285
$ws.Range("Q13").Value = 0.02401401194866667
$ws.Range("R13").Value = 0.216126107538
$ws.Range("S13").Value = 0.003725474212292008
$ws.Range("T13").Value = 0.003725474212292007

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.05732366666666666
$ws.Range("H14").Value = 0.171971
$ws.Range("I14").Value = 0.05075379138680971
$ws.Range("J14").Value = 0.05075379138680969
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.949891
$ws.Range("N14").Value = 2.849673
$ws.Range("O14").Value = 0.1664393778377885
$ws.Range("P14").Value = 0.1664393778377885
$ws.Range("Q14").Value = 0.05445123505366666
$ws.Range("R14").Value = 0.490061115483
$ws.Range("S14").Value = 0.008447429461329517
$ws.Range("T14").Value = 0.008447429461329514

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.05732366666666666
$ws.Range("H15").Value = 0.171971
$ws.Range("I15").Value = 0.05075379138680971
$ws.Range("J15").Value = 0.05075379138680969
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.834886333333333
$ws.Range("N15").Value = 8.504659
$ws.Range("O15").Value = 0.4967272219242518
$ws.Range("P15").Value = 0.4967272219242519
$ws.Range("Q15").Value = 0.1625060792098889
$ws.Range("R15").Value = 1.462554712889
$ws.Range("S15").Value = 0.02521078979769301
$ws.Range("T15").Value = 0.025210789797693

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.05732366666666666
$ws.Range("H16").Value = 0.171971
$ws.Range("I16").Value = 0.05075379138680971
$ws.Range("J16").Value = 0.05075379138680969
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.722217666666667
$ws.Range("N16").Value = 5.166653
$ws.Range("O16").Value = 0.3017660309880268
$ws.Range("P16").Value = 0.3017660309880269
$ws.Range("Q16").Value = 0.09872383145144444
$ws.Range("R16").Value = 0.888514483063
$ws.Range("S16").Value = 0.01531577018439187
$ws.Range("T16").Value = 0.01531577018439187

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.05732366666666666
$ws.Range("H17").Value = 0.171971
$ws.Range("I17").Value = 0.05075379138680971
$ws.Range("J17").Value = 0.05075379138680969
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.200134
$ws.Range("N17").Value = 0.600402
$ws.Range("O17").Value = 0.03975887129191588
$ws.Range("P17").Value = 0.03975887129191589
$ws.Range("Q17").Value = 0.01147241470466667
$ws.Range("R17").Value = 0.103251732342
$ws.Range("S17").Value = 0.001779801943395317
$ws.Range("T17").Value = 0.001779801943395317

